$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 4): No=10, Level Kode=DRTU, Level Nama=Direktur Keuangan
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "DRTU"
$ws.Range("C4").Value = "Direktur Keuangan"

# Update the active selection to match the author's final cursor position
$ws.Range("C12").Select()
